$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Ccl25/Ackr2 -> FAPs)
$ws.Range("G2").Value = 2.738607
$ws.Range("H2").Value = 8.215821
$ws.Range("I2").Value = 0.2235648590725649
$ws.Range("J2").Value = 0.223564859072565
$ws.Range("Q2").Value = 33.838427097442
$ws.Range("R2").Value = 304.545843876978
$ws.Range("S2").Value = 0.2235648590725649
$ws.Range("T2").Value = 0.223564859072565

# Row 3 (FAPs -> Ccl25/Ackr2 -> FAPs)
$ws.Range("I3").Value = 0.5119261090069511
$ws.Range("J3").Value = 0.5119261090069511
$ws.Range("S3").Value = 0.5119261090069511
$ws.Range("T3").Value = 0.5119261090069511

# Row 4 (MuSCs -> Ccl25/Ackr2 -> FAPs)
$ws.Range("I4").Value = 0.2645090319204839
$ws.Range("J4").Value = 0.2645090319204839
$ws.Range("S4").Value = 0.2645090319204839
$ws.Range("T4").Value = 0.2645090319204839
